# ----------------------------------------------------------------------------
# cryptos.xlsx refresh - GitHub Actions scheduled price/volume update
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.203.07"
$ws.Range("E2").Value = "  +1.03%  "

$ws.Range("D3").Value = "1.783.79"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.32"
$ws.Range("E5").Value = "  +0.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.547"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.81"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0691"
$ws.Range("E10").Value = "  +2.10%  "

$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("D12").Value = "2.041.70"
$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.95"
$ws.Range("E13").Value = "  -2.33%  "

$ws.Range("D14").Value = "1.774.28"
$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "34.177.04"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.624"
$ws.Range("E16").Value = "  +2.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.18"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.93"
$ws.Range("E18").Value = "  +1.98%  "

$ws.Range("D19").Value = "0.0₃0801"
$ws.Range("E19").Value = "  +3.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.33"
$ws.Range("E20").Value = "  +3.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.97"
$ws.Range("E21").Value = "  +3.75%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  +2.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -1.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.36"
$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.20"
$ws.Range("E26").Value = "  +2.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.30"
$ws.Range("E27").Value = "  +1.35%  "

$ws.Range("E28").Value = "  +1.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.40%  "

$ws.Range("E30").Value = "  +0.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0520"
$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.74"
$ws.Range("E32").Value = "  +4.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.74"
$ws.Range("E33").Value = "  +6.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.79"
$ws.Range("E34").Value = "  -1.20%  "

$ws.Range("D35").Value = "1.444.01"
$ws.Range("E35").Value = "  +3.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.654"
$ws.Range("E36").Value = "  +3.00%  "

$ws.Range("E37").Value = "  +6.25%  "

$ws.Range("E38").Value = "  +3.25%  "

$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.26"
$ws.Range("E40").Value = "  +2.41%  "

$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.924"
$ws.Range("E42").Value = "  +1.38%  "

$ws.Range("E43").Value = "  +0.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.52"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0510"
$ws.Range("E45").Value = "  +0.36%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.08"
$ws.Range("E46").Value = "  +3.75%  "

$ws.Range("E47").Value = "  +0.05%  "

$ws.Range("D48").Value = "0.0⁦0134"
$ws.Range("E48").Value = "  -3.70%  "

$ws.Range("D49").Value = "1.943.00"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.44"
$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("E51").Value = "  +0.18%  "
